$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 247, pushing existing rows 247-311 down to 248-312
$ws.Rows(247).Insert()

# Populate the new row 247 with the new weekly record
$ws.Cells.Item(247, 1).Value = 3
$ws.Cells.Item(247, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(247, 3).Value = "Coquimbo"
$ws.Cells.Item(247, 4).Value = 44642
$ws.Cells.Item(247, 5).Value = 5
$ws.Cells.Item(247, 6).Value = 100112009
$ws.Cells.Item(247, 7).Value = "Acelga"
$ws.Cells.Item(247, 8).Value = "Sin especificar"
$ws.Cells.Item(247, 9).Value = "Primera"
$ws.Cells.Item(247, 10).Value = 185
$ws.Cells.Item(247, 11).Value = 4000
$ws.Cells.Item(247, 12).Value = 4500
$ws.Cells.Item(247, 13).Value = 4243
$ws.Cells.Item(247, 14).Value = "$/docena de atados (6 kilos)"
$ws.Cells.Item(247, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(247, 16).Value = 707
$ws.Cells.Item(247, 17).Value = 6
$ws.Cells.Item(247, 18).Value = "Hortaliza"
